# Incorporo nuevos datos hasta diciembre de 2025
# Update column H (year 2025) values for rows 2-18 on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 10.3887249116102
    3  = 11.38948747670552
    4  = 11.92926297162149
    5  = 10.95093022876864
    6  = 9.849122963011057
    7  = 10.64535675891659
    8  = 11.01244336109707
    9  = 10.40394883131215
    10 = 12.96039803917507
    11 = 10.54511212492152
    12 = 10.2869594890205
    13 = 11.01746589743031
    14 = 13.88131005354286
    15 = 10.6675283336404
    16 = 13.4013689896905
    17 = 14.60795920176889
    18 = 11.1941193388301
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $newValues[$row]
}
